# --- Rename the "Requested quantity" columns on the existing sheets ---
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$ws2 = $wb.Worksheets.Item(2)   # "Monthly Trend"

$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Match the page margins used by the other sheets (0.75in/1in/0.5in)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Reuse the existing bold/bordered header style from "Weekly Quantity"!A1:B1
$ws1.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Reuse the existing date-formatted style from "Weekly Quantity"!A2
$ws1.Range("A2").Copy()
$newSheet.Range("A2:A41").PasteSpecial(-4122)

# --- Header row ---
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# --- PO forecast data (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$newSheet.Cells.Item(2,1).Value = 45319.99999999999
$newSheet.Cells.Item(2,2).Value = 165
$newSheet.Cells.Item(2,3).Value = -1.606381319722833
$newSheet.Cells.Item(2,4).Value = 337.0673299896298
$newSheet.Cells.Item(3,1).Value = 45333.99999999999
$newSheet.Cells.Item(3,2).Value = 167
$newSheet.Cells.Item(3,3).Value = -6.712938071644828
$newSheet.Cells.Item(3,4).Value = 331.4598667913543
$newSheet.Cells.Item(4,1).Value = 45340.99999999999
$newSheet.Cells.Item(4,2).Value = 168
$newSheet.Cells.Item(4,3).Value = -18.59711561298863
$newSheet.Cells.Item(4,4).Value = 317.2355249908561
$newSheet.Cells.Item(5,1).Value = 45347.99999999999
$newSheet.Cells.Item(5,2).Value = 169
$newSheet.Cells.Item(5,3).Value = 4.253688031603632
$newSheet.Cells.Item(5,4).Value = 338.795225870944
$newSheet.Cells.Item(6,1).Value = 45354.99999999999
$newSheet.Cells.Item(6,2).Value = 169
$newSheet.Cells.Item(6,3).Value = 0.9447191887366607
$newSheet.Cells.Item(6,4).Value = 345.0702410872705
$newSheet.Cells.Item(7,1).Value = 45368.99999999999
$newSheet.Cells.Item(7,2).Value = 171
$newSheet.Cells.Item(7,3).Value = 7.07067195624935
$newSheet.Cells.Item(7,4).Value = 341.8580619497407
$newSheet.Cells.Item(8,1).Value = 45389.99999999999
$newSheet.Cells.Item(8,2).Value = 174
$newSheet.Cells.Item(8,3).Value = 15.36763777141555
$newSheet.Cells.Item(8,4).Value = 334.2956518526397
$newSheet.Cells.Item(9,1).Value = 45403.99999999999
$newSheet.Cells.Item(9,2).Value = 176
$newSheet.Cells.Item(9,3).Value = 28.83934149204891
$newSheet.Cells.Item(9,4).Value = 339.1300451078206
$newSheet.Cells.Item(10,1).Value = 45410.99999999999
$newSheet.Cells.Item(10,2).Value = 177
$newSheet.Cells.Item(10,3).Value = 14.15627950812373
$newSheet.Cells.Item(10,4).Value = 344.0128776454534
$newSheet.Cells.Item(11,1).Value = 45417.99999999999
$newSheet.Cells.Item(11,2).Value = 178
$newSheet.Cells.Item(11,3).Value = 7.650786922199567
$newSheet.Cells.Item(11,4).Value = 340.1419981208458
$newSheet.Cells.Item(12,1).Value = 45424.99999999999
$newSheet.Cells.Item(12,2).Value = 179
$newSheet.Cells.Item(12,3).Value = 23.12157471915499
$newSheet.Cells.Item(12,4).Value = 352.1967238018545
$newSheet.Cells.Item(13,1).Value = 45431.99999999999
$newSheet.Cells.Item(13,2).Value = 180
$newSheet.Cells.Item(13,3).Value = 17.90237606135699
$newSheet.Cells.Item(13,4).Value = 345.6707974426129
$newSheet.Cells.Item(14,1).Value = 45445.99999999999
$newSheet.Cells.Item(14,2).Value = 182
$newSheet.Cells.Item(14,3).Value = 23.35633261925196
$newSheet.Cells.Item(14,4).Value = 350.2901881488587
$newSheet.Cells.Item(15,1).Value = 45459.99999999999
$newSheet.Cells.Item(15,2).Value = 184
$newSheet.Cells.Item(15,3).Value = 29.7634276480652
$newSheet.Cells.Item(15,4).Value = 355.9123933135411
$newSheet.Cells.Item(16,1).Value = 45466.99999999999
$newSheet.Cells.Item(16,2).Value = 185
$newSheet.Cells.Item(16,3).Value = 10.54965677935609
$newSheet.Cells.Item(16,4).Value = 338.3878635535044
$newSheet.Cells.Item(17,1).Value = 45473.99999999999
$newSheet.Cells.Item(17,2).Value = 186
$newSheet.Cells.Item(17,3).Value = 19.67229509356508
$newSheet.Cells.Item(17,4).Value = 357.3593766654562
$newSheet.Cells.Item(18,1).Value = 45480.99999999999
$newSheet.Cells.Item(18,2).Value = 187
$newSheet.Cells.Item(18,3).Value = 29.66394207955831
$newSheet.Cells.Item(18,4).Value = 351.5746663957887
$newSheet.Cells.Item(19,1).Value = 45487.99999999999
$newSheet.Cells.Item(19,2).Value = 188
$newSheet.Cells.Item(19,3).Value = 21.41384853545142
$newSheet.Cells.Item(19,4).Value = 350.2894133646889
$newSheet.Cells.Item(20,1).Value = 45494.99999999999
$newSheet.Cells.Item(20,2).Value = 189
$newSheet.Cells.Item(20,3).Value = 16.62238596241329
$newSheet.Cells.Item(20,4).Value = 342.7558593729057
$newSheet.Cells.Item(21,1).Value = 45501.99999999999
$newSheet.Cells.Item(21,2).Value = 190
$newSheet.Cells.Item(21,3).Value = 44.89875697346812
$newSheet.Cells.Item(21,4).Value = 347.9384235034392
$newSheet.Cells.Item(22,1).Value = 45515.99999999999
$newSheet.Cells.Item(22,2).Value = 192
$newSheet.Cells.Item(22,3).Value = 9.758002266493628
$newSheet.Cells.Item(22,4).Value = 352.2274920611293
$newSheet.Cells.Item(23,1).Value = 45522.99999999999
$newSheet.Cells.Item(23,2).Value = 193
$newSheet.Cells.Item(23,3).Value = 46.01230205442248
$newSheet.Cells.Item(23,4).Value = 360.7127360680165
$newSheet.Cells.Item(24,1).Value = 45529.99999999999
$newSheet.Cells.Item(24,2).Value = 194
$newSheet.Cells.Item(24,3).Value = 27.18700386363907
$newSheet.Cells.Item(24,4).Value = 352.0183765300732
$newSheet.Cells.Item(25,1).Value = 45536.99999999999
$newSheet.Cells.Item(25,2).Value = 195
$newSheet.Cells.Item(25,3).Value = 38.86295933249171
$newSheet.Cells.Item(25,4).Value = 358.7499305056629
$newSheet.Cells.Item(26,1).Value = 45543.99999999999
$newSheet.Cells.Item(26,2).Value = 196
$newSheet.Cells.Item(26,3).Value = 34.63044201313476
$newSheet.Cells.Item(26,4).Value = 371.268646463647
$newSheet.Cells.Item(27,1).Value = 45550.99999999999
$newSheet.Cells.Item(27,2).Value = 197
$newSheet.Cells.Item(27,3).Value = 37.71672005583293
$newSheet.Cells.Item(27,4).Value = 355.8228365081702
$newSheet.Cells.Item(28,1).Value = 45557.99999999999
$newSheet.Cells.Item(28,2).Value = 198
$newSheet.Cells.Item(28,3).Value = 41.59574434820199
$newSheet.Cells.Item(28,4).Value = 370.5636102604956
$newSheet.Cells.Item(29,1).Value = 45564.99999999999
$newSheet.Cells.Item(29,2).Value = 199
$newSheet.Cells.Item(29,3).Value = 40.91553618282772
$newSheet.Cells.Item(29,4).Value = 358.4840820603981
$newSheet.Cells.Item(30,1).Value = 45585.99999999999
$newSheet.Cells.Item(30,2).Value = 202
$newSheet.Cells.Item(30,3).Value = 21.46800038496522
$newSheet.Cells.Item(30,4).Value = 359.9828751603424
$newSheet.Cells.Item(31,1).Value = 45592.99999999999
$newSheet.Cells.Item(31,2).Value = 203
$newSheet.Cells.Item(31,3).Value = 34.46477666397477
$newSheet.Cells.Item(31,4).Value = 371.1163769693252
$newSheet.Cells.Item(32,1).Value = 45599.99999999999
$newSheet.Cells.Item(32,2).Value = 204
$newSheet.Cells.Item(32,3).Value = 35.3882502981881
$newSheet.Cells.Item(32,4).Value = 368.5997896956047
$newSheet.Cells.Item(33,1).Value = 45613.99999999999
$newSheet.Cells.Item(33,2).Value = 206
$newSheet.Cells.Item(33,3).Value = 26.66937960022772
$newSheet.Cells.Item(33,4).Value = 374.6567487979323
$newSheet.Cells.Item(34,1).Value = 45620.99999999999
$newSheet.Cells.Item(34,2).Value = 207
$newSheet.Cells.Item(34,3).Value = 44.41977513414781
$newSheet.Cells.Item(34,4).Value = 379.5688662879995
$newSheet.Cells.Item(35,1).Value = 45627.99999999999
$newSheet.Cells.Item(35,2).Value = 208
$newSheet.Cells.Item(35,3).Value = 40.93056682611897
$newSheet.Cells.Item(35,4).Value = 369.4717433641876
$newSheet.Cells.Item(36,1).Value = 45634.99999999999
$newSheet.Cells.Item(36,2).Value = 208
$newSheet.Cells.Item(36,3).Value = 33.36008160159775
$newSheet.Cells.Item(36,4).Value = 367.0005635588489
$newSheet.Cells.Item(37,1).Value = 45641.99999999999
$newSheet.Cells.Item(37,2).Value = 209
$newSheet.Cells.Item(37,3).Value = 50.1476464091351
$newSheet.Cells.Item(37,4).Value = 373.9499805294157
$newSheet.Cells.Item(38,1).Value = 45648.99999999999
$newSheet.Cells.Item(38,2).Value = 210
$newSheet.Cells.Item(38,3).Value = 48.61189321582133
$newSheet.Cells.Item(38,4).Value = 376.8917370998279
$newSheet.Cells.Item(39,1).Value = 45655.99999999999
$newSheet.Cells.Item(39,2).Value = 211
$newSheet.Cells.Item(39,3).Value = 48.44340392930093
$newSheet.Cells.Item(39,4).Value = 383.8145873209214
$newSheet.Cells.Item(40,1).Value = 45662.99999999999
$newSheet.Cells.Item(40,2).Value = 212
$newSheet.Cells.Item(40,3).Value = 37.88626929400352
$newSheet.Cells.Item(40,4).Value = 374.4653421657113
$newSheet.Cells.Item(41,1).Value = 45669.99999999999
$newSheet.Cells.Item(41,2).Value = 213
$newSheet.Cells.Item(41,3).Value = 44.98421887023439
$newSheet.Cells.Item(41,4).Value = 377.3193373133338
